$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (shifts old J:L -> K:M).
# The new column inherits column I's width/format, matching Excel's
# default "insert column" behaviour.
$ws.Columns("J").Insert()
$ws.Columns("J").ColumnWidth = $ws.Columns("I").ColumnWidth

# Split the old single "ExpectedSourceTemplateFile" column into separate
# Excel/Word expected-template-file columns.
$ws.Range("I1").Value = "ExpectedSourceTemplateFile_Excel"
$ws.Range("J1").Value = "ExpectedSourceTemplateFile_Word"

# Populate the new "Word" template paths per report category row.
$ws.Range("J2").Value = "\Testdata\Templates\SLRReport_SourceData\Expected_word_Data\Clinical.xlsx"
$ws.Range("J3").Value = "\Testdata\Templates\SLRReport_SourceData\Expected_word_Data\Economic.xlsx"
$ws.Range("J4").Value = "\Testdata\Templates\SLRReport_SourceData\Expected_word_Data\QOL.xlsx"
$ws.Range("J5").Value = "\Testdata\Templates\SLRReport_SourceData\Expected_word_Data\RWE.xlsx"

# Update the sheet view: move the active selection (the sheet also
# scrolls back so column A is visible again).
$null = $ws.Range("D12").Select()
